$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.091.26"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -1.12%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.677.03"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -0.59%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D5').Value = "'210.96"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -3.53%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'0.5288"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -4.50%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'1.004"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  -0.44%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.2679"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  -1.17%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.06316"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -2.73%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'21.29"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -3.84%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.07575"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -0.10%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'1.680.36"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -0.43%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'4.509"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -1.07%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'0.5683"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -2.39%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'0.000008142"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -4.00%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('E16').Value = "'  +0.40%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'26.131.54"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -1.08%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'1.004"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -0.39%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'4.864"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -1.68%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'  -3.12%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'189.45"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -0.95%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'6.208"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -0.64%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = "'  -0.42%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'148.60"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -0.43%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'0.1259"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -4.80%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'7.655"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -3.25%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'16.05"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +1.39%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'0.06374"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  +0.54%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'1.348"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -3.46%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'1.286"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -3.22%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'3.541"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -1.37%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'3.540"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  -1.21%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'1.673"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -0.22%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'  -3.07%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'0.6067"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -2.87%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'2.418"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +0.35%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'2.721"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +0.21%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'6.157"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -1.37%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'0.01615"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -1.38%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'1.095.64"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  -1.89%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'0.8714"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -0.71%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = "'  -1.01%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'100.07"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -0.62%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'1.828.47"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -0.44%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'0.00000000109"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -1.90%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'57.10"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -0.67%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'1.007"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +0.12%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('B48').Value = "'EnergySwap"
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').Value = "'8.012"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -2.61%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').Value = "'Cronos"
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').Value = "'0.05253"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -0.62%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'0.4265"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -0.79%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'5.963"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -2.10%  "
$ws.Range('E51').Style = 'Normal'
